$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# SafeArray of points?
$pts = @(@(10,10), @(50,50), @(90,10))
try {
    $newShape = $ws.Shapes.AddPolyline($pts)
    Write-Output "New shape name: $($newShape.Name())"
} catch {
    Write-Output "ERROR: $_"
}
